# How to get and update the data from excel based on filter search criteria
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header/sno column label
$ws.Range("A1").Value = "350"

# Update price cell to match filter result (text "350" instead of numeric 299)
$ws.Range("D2").Value = "350"

# Fix fruit name typo/replacement
$ws.Range("B3").Value = "Apple"

# Fix fruit name typo
$ws.Range("B5").Value = "Rebulic"
